$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5 (scenario 3 - NavigateToComment): description now refers to the blog pages.
$ws.Range("C5").Value = "This is to test whether users are able to successfully navigate to the blog pages to leave comments"

# Row 6 - new test case: NavigateToBlogPages
$ws.Range("B6").Value = "test_<NavigateToBlogPages>"
$ws.Range("C6").Value = "This is to test whether users are able to navigate to each blog pages"
$ws.Range("D6").Value = "NIL"
$ws.Range("E6").ClearFormats()
$ws.Range("E6").Value = "Navigated to each Blog pages"
$ws.Range("F6").Value = "Navigated to each Blog pages"

# Row 7 - new test case: NavigateToProjectPages
$ws.Range("B7").Value = "test_<NavigateToProjectPages>"
$ws.Range("C7").Value = "This is to test whether users are able to navigate to each project pages"
$ws.Range("D7").Value = "NIL"
$ws.Range("E7").Value = "Navigated to each Project pages"
$ws.Range("F7").Value = "Navigated to each Project pages"

# Row 8 - new test case: ViewResume
$ws.Range("B8").Value = "test_<ViewResume>"
$ws.Range("C8").Value = "This is to test whether users are able to view a created Resume"
$ws.Range("D8").Value = "NIL"
$ws.Range("E8").Value = "Resume page is shown"
$ws.Range("F8").Value = "Resume page is shown"

# Row 9 - new test case: SelfIntro
$ws.Range("B9").Value = "test_<SelfIntro>"
$ws.Range("C9").Value = "This is to test whether users are able to view a created SelfIntro"
$ws.Range("D9").Value = "NIL"
$ws.Range("E9").Value = "Self Intro page is shown"
$ws.Range("F9").Value = "Self Intro page is shown"

# Justification column for every test case that has already passed.
$ws.Range("G3").Value = "Based on the given source code, this function has already been implemented"
$ws.Range("G4").Value = "Based on the given source code, this function has already been implemented"
$ws.Range("G5").Value = "Based on the given source code, this function has already been implemented"
$ws.Range("G6").Value = "Based on the given source code, this function has already been implemented"
$ws.Range("G7").Value = "Based on the given source code, this function has already been implemented"
$ws.Range("G8").Value = "Based on the given source code, I created a blog post as a Resume in the /admin page"
$ws.Range("G9").Value = "Based on the given source code, I created a blog post as a Resume in the /admin page"

# Update the active selection to match the saved view.
$ws.Range("B10").Select()
